$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.61"
$ws.Range("E2").Value = "'2.83%"
$ws.Range("G2").Value = "'12"

$ws.Range("D3").Value = "'39.45"
$ws.Range("E3").Value = "'2.25%"
$ws.Range("G3").Value = "'12"

$ws.Range("D4").Value = "'5.136"
$ws.Range("E4").Value = "'0.63%"
$ws.Range("G4").Value = "'12"

$ws.Range("D5").Value = "'0.08190"
$ws.Range("E5").Value = "'0.98%"
$ws.Range("G5").Value = "'12"

$ws.Range("D6").Value = "'1.972"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("G6").Value = "'12"

$ws.Range("D7").Value = "'8.230"
$ws.Range("E7").Value = "'3.53%"
$ws.Range("G7").Value = "'12"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.245"
$ws.Range("E8").Value = "'1.27%"
$ws.Range("G8").Value = "'12"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9279"
$ws.Range("E9").Value = "'-0.32%"
$ws.Range("G9").Value = "'12"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1414"
$ws.Range("E10").Value = "'-3.02%"
$ws.Range("G10").Value = "'12"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1992"
$ws.Range("E11").Value = "'1.88%"
$ws.Range("G11").Value = "'12"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09077"
$ws.Range("E12").Value = "'-0.30%"
$ws.Range("G12").Value = "'12"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03498"
$ws.Range("E13").Value = "'-0.31%"
$ws.Range("G13").Value = "'12"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09818"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("G14").Value = "'12"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001398"
$ws.Range("E15").Value = "'-0.21%"
$ws.Range("G15").Value = "'12"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005870"
$ws.Range("E16").Value = "'-3.27%"
$ws.Range("G16").Value = "'12"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.649"
$ws.Range("E17").Value = "'-2.02%"
$ws.Range("G17").Value = "'12"

$ws.Range("D18").Value = "'3.163"
$ws.Range("E18").Value = "'-8.22%"
$ws.Range("G18").Value = "'12"

$ws.Range("D19").Value = "'0.3466"
$ws.Range("E19").Value = "'0.11%"
$ws.Range("G19").Value = "'12"

$ws.Range("E20").Value = "'0.79%"
$ws.Range("G20").Value = "'12"

$ws.Range("D21").Value = "'4.842"
$ws.Range("E21").Value = "'0.58%"
$ws.Range("G21").Value = "'12"

$ws.Range("D22").Value = "'0.2447"
$ws.Range("E22").Value = "'-0.21%"
$ws.Range("G22").Value = "'12"

$ws.Range("D23").Value = "'0.04367"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("G23").Value = "'12"

$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("G24").Value = "'12"

$ws.Range("D25").Value = "'0.004784"
$ws.Range("E25").Value = "'-1.05%"
$ws.Range("G25").Value = "'12"

$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.25%"
$ws.Range("G26").Value = "'12"

$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'-10.16%"
$ws.Range("G27").Value = "'12"

$ws.Range("G28").Value = "'12"

$ws.Range("G29").Value = "'12"

$ws.Range("G30").Value = "'12"

$ws.Range("G31").Value = "'12"

$ws.Range("G32").Value = "'12"

$ws.Range("G33").Value = "'12"

$ws.Range("G34").Value = "'12"

$ws.Range("G35").Value = "'12"

$ws.Range("G36").Value = "'12"

$ws.Range("G37").Value = "'12"

$ws.Range("G38").Value = "'12"

$ws.Range("D39").Value = "'0.02207"
$ws.Range("E39").Value = "'5.28%"
$ws.Range("G39").Value = "'12"

$ws.Range("D40").Value = "'0.05185"
$ws.Range("E40").Value = "'1.34%"
$ws.Range("G40").Value = "'12"

$ws.Range("E41").Value = "'1.40%"
$ws.Range("G41").Value = "'12"

$ws.Range("D42").Value = "'0.009778"
$ws.Range("E42").Value = "'-3.61%"
$ws.Range("G42").Value = "'12"

$ws.Range("D43").Value = "'0.1374"
$ws.Range("E43").Value = "'1.05%"
$ws.Range("G43").Value = "'12"

$ws.Range("D44").Value = "'0.002128"
$ws.Range("E44").Value = "'-0.25%"
$ws.Range("G44").Value = "'12"

$ws.Range("D45").Value = "'0.009194"
$ws.Range("E45").Value = "'-0.66%"
$ws.Range("G45").Value = "'12"

$ws.Range("D46").Value = "'0.00006390"
$ws.Range("E46").Value = "'2.33%"
$ws.Range("G46").Value = "'12"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("G47").Value = "'12"

$ws.Range("D48").Value = "'0.002763"
$ws.Range("E48").Value = "'-8.77%"
$ws.Range("G48").Value = "'12"

$ws.Range("D49").Value = "'0.001199"
$ws.Range("E49").Value = "'-25.14%"
$ws.Range("G49").Value = "'12"

$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("G50").Value = "'12"

$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.26%"
$ws.Range("G51").Value = "'12"

